# Auto-generated edit script applying numeric cell updates
# across sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1441.5714
$ws.Range("I31").Value = 1038.2
$ws.Range("K31").Value = 3114.6
$ws.Range("M31").Value = -2884.6
$ws.Range("H32").Value = 1183.3334
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 1575
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 1575
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -2227
$ws.Range("H100").Value = 2228
$ws.Range("I100").Value = 2358
$ws.Range("J100").Value = 1903
$ws.Range("K100").Value = 2358
$ws.Range("L100").Value = 1903
$ws.Range("M100").Value = -1817
$ws.Range("N100").Value = -2985
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3573.5454
$ws.Range("I102").Value = 3602
$ws.Range("K102").Value = 3602
$ws.Range("M102").Value = -1980
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1705.6666
$ws.Range("I99").Value = 1932.9231
$ws.Range("J99").Value = 1336.375
$ws.Range("K99").Value = 1932.9231
$ws.Range("L99").Value = 1336.375
$ws.Range("M99").Value = -434.9231
$ws.Range("N99").Value = -4332.375
$ws.Range("H105").Value = 5491.8184
$ws.Range("I105").Value = 5491.8184
$ws.Range("K105").Value = 5491.8184
$ws.Range("M105").Value = -3744.8184
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2641.2646
$ws.Range("I31").Value = 2754.125
$ws.Range("J31").Value = 2370.4
$ws.Range("K31").Value = 2754.125
$ws.Range("L31").Value = 2370.4
$ws.Range("M31").Value = -2459.125
$ws.Range("N31").Value = -2960.4
$ws.Range("H34").Value = 2641.2646
$ws.Range("I34").Value = 2754.125
$ws.Range("J34").Value = 2370.4
$ws.Range("K34").Value = 2754.125
$ws.Range("L34").Value = 2370.4
$ws.Range("M34").Value = -2552.125
$ws.Range("N34").Value = -2774.4
$ws.Range("H58").Value = 3582.3428
$ws.Range("I58").Value = 544.6923
$ws.Range("J58").Value = 12357.777
$ws.Range("K58").Value = 544.6923
$ws.Range("L58").Value = 12357.777
$ws.Range("M58").Value = -341.6923
$ws.Range("N58").Value = -12763.777
$ws.Range("H98").Value = 44666.668
$ws.Range("J98").Value = 44666.668
$ws.Range("L98").Value = 44666.668
$ws.Range("N98").Value = -49158.668
$ws.Range("H132").Value = 2391.8
$ws.Range("I132").Value = 1206.5333
$ws.Range("J132").Value = 5947.6
$ws.Range("K132").Value = 3619.5999
$ws.Range("L132").Value = 17842.8
$ws.Range("M132").Value = -1089.5999
$ws.Range("N132").Value = -22902.8
$ws.Range("H136").Value = 3582.3428
$ws.Range("I136").Value = 544.6923
$ws.Range("J136").Value = 12357.777
$ws.Range("K136").Value = 1634.0769
$ws.Range("L136").Value = 37073.331
$ws.Range("M136").Value = 915.9231
$ws.Range("N136").Value = -42173.331
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("K16").Value = 300
$ws.Range("M16").Value = -127
$ws.Range("H45").Value = 923.25
$ws.Range("I45").Value = 540
$ws.Range("K45").Value = 1620
$ws.Range("M45").Value = -1088
$ws.Range("H49").Value = 893.75
$ws.Range("I49").Value = 893.75
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2681.25
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2525.25
$ws.Range("N49").ClearContents()
$ws.Range("H58").Value = 2995
$ws.Range("I58").Value = 1005
$ws.Range("J58").Value = 3990
$ws.Range("K58").Value = 3015
$ws.Range("L58").Value = 11970
$ws.Range("N58").Value = -12226
$ws.Range("M58").Value = -2887
$ws.Range("H63").Value = 4020
$ws.Range("I63").Value = 800
$ws.Range("J63").Value = 4377.778
$ws.Range("K63").Value = 2400
$ws.Range("L63").Value = 13133.334
$ws.Range("M63").Value = -1651
$ws.Range("N63").Value = -14631.334
$ws.Range("H66").Value = 4020
$ws.Range("I66").Value = 800
$ws.Range("J66").Value = 4377.778
$ws.Range("K66").Value = 7200
$ws.Range("L66").Value = 39400.002
$ws.Range("M66").Value = -3456
$ws.Range("N66").Value = -46888.002
$ws.Range("H69").Value = 866.6667
$ws.Range("J69").Value = 866.6667
$ws.Range("L69").Value = 2600.0001
$ws.Range("N69").Value = -4222.0001
$ws.Range("H72").Value = 866.6667
$ws.Range("J72").Value = 866.6667
$ws.Range("L72").Value = 7800.0003
$ws.Range("N72").Value = -15912.0003
$ws.Range("H88").Value = 8333.333000000001
$ws.Range("J88").Value = 8333.333000000001
$ws.Range("L88").Value = 24999.999
$ws.Range("N88").Value = -25855.999
$ws.Range("H91").Value = 8333.333000000001
$ws.Range("J91").Value = 8333.333000000001
$ws.Range("L91").Value = 24999.999
$ws.Range("N91").Value = -27963.999
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 1112322.1
$ws.Range("I122").Value = 466.66666
$ws.Range("J122").Value = 1668249.9
$ws.Range("K122").Value = 4199.99994
$ws.Range("L122").Value = 15014249.1
$ws.Range("M122").Value = -1749.99994
$ws.Range("N122").Value = -15019149.1
$ws.Range("H133").Value = 10000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -40120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H93").Value = 1363.4286
$ws.Range("I93").Value = 1369.1
$ws.Range("J93").Value = 1349.25
$ws.Range("K93").Value = 1369.1
$ws.Range("L93").Value = 1349.25
$ws.Range("M93").Value = -121.0999999999999
$ws.Range("N93").Value = -3845.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1338.975
$ws.Range("I136").Value = 1297.9242
$ws.Range("J136").Value = 1532.5
$ws.Range("K136").Value = 3893.7726
$ws.Range("L136").Value = 4597.5
$ws.Range("M136").Value = -1343.7726
$ws.Range("N136").Value = -9697.5
